$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nRows = 24
$block1 = New-Object 'double[,]' 24,5
$block1[0,0] = 1.02
$block1[0,1] = 1.066139318179095
$block1[0,2] = 1.067857451690567
$block1[0,3] = 0.992614727750844
$block1[0,4] = 1.076761132788726
$block1[1,0] = 1.02
$block1[1,1] = 1.067404269981731
$block1[1,2] = 1.068877336879326
$block1[1,3] = 0.9936372048519299
$block1[1,4] = 1.077977673300468
$block1[2,0] = 1.02
$block1[2,1] = 1.068221776081627
$block1[2,2] = 1.069536277956044
$block1[2,3] = 0.9942998659930998
$block1[2,4] = 1.078764247416554
$block1[3,0] = 1.02
$block1[3,1] = 1.0685652192807
$block1[3,2] = 1.069813061425091
$block1[3,3] = 0.994578699834602
$block1[3,4] = 1.079094780219495
$block1[4,0] = 1.02
$block1[4,1] = 1.068622871148538
$block1[4,2] = 1.069859520853329
$block1[4,3] = 0.994625531979634
$block1[4,4] = 1.079150269830519
$block1[5,0] = 1.02
$block1[5,1] = 1.068226366110747
$block1[5,2] = 1.06953997727221
$block1[5,3] = 0.994303590798249
$block1[5,4] = 1.078768664570239
$block1[6,0] = 1.02
$block1[6,1] = 1.066567024259258
$block1[6,2] = 1.068202333643187
$block1[6,3] = 0.9929600610674297
$block1[6,4] = 1.077172396140735
$block1[7,0] = 1.02
$block1[7,1] = 1.063635214335322
$block1[7,2] = 1.065837525861173
$block1[7,3] = 0.9906006454969559
$block1[7,4] = 1.074354776286639
$block1[8,0] = 1.02
$block1[8,1] = 1.061675177355689
$block1[8,2] = 1.064255655378965
$block1[8,3] = 0.989033133672735
$block1[8,4] = 1.072472961640151
$block1[9,0] = 1.02
$block1[9,1] = 1.060825104586827
$block1[9,2] = 1.063569389425585
$block1[9,3] = 0.988355674866747
$block1[9,4] = 1.071657267479959
$block1[10,0] = 1.02
$block1[10,1] = 1.060509140296651
$block1[10,2] = 1.063314280524437
$block1[10,3] = 0.9881042295826724
$block1[10,4] = 1.071354150125639
$block1[11,0] = 1.02
$block1[11,1] = 1.060576925258641
$block1[11,2] = 1.063369011311251
$block1[11,3] = 0.9881581567098651
$block1[11,4] = 1.07141917586887
$block1[12,0] = 1.02
$block1[12,1] = 1.06079899117822
$block1[12,2] = 1.063548306124768
$block1[12,3] = 0.9883348863814464
$block1[12,4] = 1.071632214423239
$block1[13,0] = 1.02
$block1[13,1] = 1.060935785435494
$block1[13,2] = 1.063658749094279
$block1[13,3] = 0.9884438009545853
$block1[13,4] = 1.071763456878105
$block1[14,0] = 1.02
$block1[14,1] = 1.061731564766624
$block1[14,2] = 1.064301172825427
$block1[14,3] = 0.9890781214508737
$block1[14,4] = 1.072527078231646
$block1[15,0] = 1.02
$block1[15,1] = 1.062230367801363
$block1[15,2] = 1.064703796587819
$block1[15,3] = 0.9894763578477731
$block1[15,4] = 1.073005846006348
$block1[16,0] = 1.02
$block1[16,1] = 1.062521179890472
$block1[16,2] = 1.0649385144513
$block1[16,3] = 0.9897087662937551
$block1[16,4] = 1.073285020804198
$block1[17,0] = 1.02
$block1[17,1] = 1.062620317165764
$block1[17,2] = 1.065018525935803
$block1[17,3] = 0.9897880325774039
$block1[17,4] = 1.073380198337393
$block1[18,0] = 1.02
$block1[18,1] = 1.062176864569555
$block1[18,2] = 1.064660611899953
$block1[18,3] = 0.9894336180360677
$block1[18,4] = 1.072954487299935
$block1[19,0] = 1.02
$block1[19,1] = 1.060733604116985
$block1[19,2] = 1.063495513793394
$block1[19,3] = 0.9882828385668249
$block1[19,4] = 1.071569483558035
$block1[20,0] = 1.02
$block1[20,1] = 1.059824955129217
$block1[20,2] = 1.062761816527222
$block1[20,3] = 0.9875604150241495
$block1[20,4] = 1.070697909918567
$block1[21,0] = 1.02
$block1[21,1] = 1.060306763722035
$block1[21,2] = 1.063150873712882
$block1[21,3] = 0.9879432794636464
$block1[21,4] = 1.071160021493776
$block1[22,0] = 1.02
$block1[22,1] = 1.062201040782423
$block1[22,2] = 1.064680125590045
$block1[22,3] = 0.9894529299347244
$block1[22,4] = 1.072977694342971
$block1[23,0] = 1.02
$block1[23,1] = 1.064394109625145
$block1[23,2] = 1.066449814176102
$block1[23,3] = 0.9912096547607051
$block1[23,4] = 1.075083784807246

$ws.Range("B2:F25").Value = $block1

$block2 = New-Object 'double[,]' 24,6
$block2[0,0] = 1.05649201043406
$block2[0,1] = 1.071090847343355
$block2[0,2] = 1.070564189899962
$block2[0,3] = 0.9955398523335997
$block2[0,4] = 1.079444214072221
$block2[0,5] = 1.072611919202419
$block2[1,0] = 1.056943145109938
$block2[1,1] = 1.072010181804702
$block2[1,2] = 1.07139943376161
$block2[1,3] = 0.9963617723202687
$block2[1,4] = 1.080477344915502
$block2[1,5] = 1.073532559224151
$block2[2,0] = 1.057233199754469
$block2[2,1] = 1.072603516855246
$block2[2,2] = 1.071938289960432
$block2[2,3] = 0.9968940712668347
$block2[2,4] = 1.081144668170764
$block2[2,5] = 1.074126736878523
$block2[3,0] = 1.057354694950235
$block2[3,1] = 1.072852589760375
$block2[3,2] = 1.072164443751108
$block2[3,3] = 0.9971179600053012
$block2[3,4] = 1.081424930617541
$block2[3,5] = 1.074376163495747
$block2[4,0] = 1.057375068546347
$block2[4,1] = 1.072894388857847
$block2[4,2] = 1.072202393670645
$block2[4,3] = 0.9971555583673455
$block2[4,4] = 1.081471971556979
$block2[4,5] = 1.074418021952732
$block2[5,0] = 1.057234824921188
$block2[5,1] = 1.072606846411243
$block2[5,2] = 1.071941313334095
$block2[5,3] = 0.9968970624462089
$block2[5,4] = 1.08114841414925
$block2[5,5] = 1.074130071162872
$block2[6,0] = 1.056644859426179
$block2[6,1] = 1.071401860255377
$block2[6,2] = 1.070846797636392
$block2[6,3] = 0.9958175282591056
$block2[6,4] = 1.079793611679747
$block2[6,5] = 1.07292337378845
$block2[7,0] = 1.055590957755106
$block2[7,1] = 1.069266644406517
$block2[7,2] = 1.068905746302002
$block2[7,3] = 0.9939188001724441
$block2[7,4] = 1.07739711091143
$block2[7,5] = 1.070785125688172
$block2[8,0] = 1.054878646013741
$block2[8,1] = 1.067835021107607
$block2[8,2] = 1.067603247158138
$block2[8,3] = 0.9926553831429383
$block2[8,4] = 1.075793115426424
$block2[8,5] = 1.069351469319969
$block2[9,0] = 1.054567883162799
$block2[9,1] = 1.067213143767433
$block2[9,2] = 1.067037210040499
$block2[9,3] = 0.9921088820399291
$block2[9,4] = 1.075097029087997
$block2[9,5] = 1.068728708842641
$block2[10,0] = 1.054452100510475
$block2[10,1] = 1.066981850871436
$block2[10,2] = 1.066826648282805
$block2[10,3] = 0.9919059725120875
$block2[10,4] = 1.074838236067017
$block2[10,5] = 1.068497087484202
$block2[11,0] = 1.054476952218225
$block2[11,1] = 1.067031477608935
$block2[11,2] = 1.066871828586547
$block2[11,3] = 0.9919494934313052
$block2[11,4] = 1.074893758761934
$block2[11,5] = 1.068546784697361
$block2[12,0] = 1.05455831971379
$block2[12,1] = 1.067194031162451
$block2[12,2] = 1.067019811288319
$block2[12,3] = 0.9920921077337197
$block2[12,4] = 1.075075642000898
$block2[12,5] = 1.068709569095569
$block2[13,0] = 1.054608406281202
$block2[13,1] = 1.067294145926511
$block2[13,2] = 1.067110947100307
$block2[13,3] = 0.9921799884222134
$block2[13,4] = 1.075187674953024
$block2[13,5] = 1.068809826034076
$block2[14,0] = 1.054899221131218
$block2[14,1] = 1.067876251171274
$block2[14,2] = 1.06764076981322
$block2[14,3] = 0.9926916645766087
$block2[14,4] = 1.075839279533515
$block2[14,5] = 1.069392757935055
$block2[15,0] = 1.055081017242803
$block2[15,1] = 1.068240859364526
$block2[15,2] = 1.067972563667151
$block2[15,3] = 0.9930127773692701
$block2[15,4] = 1.07624759769115
$block2[15,5] = 1.069757883913763
$block2[16,0] = 1.055186831427946
$block2[16,1] = 1.068453338816304
$block2[16,2] = 1.068165896160329
$block2[16,3] = 0.9932001317071766
$block2[16,4] = 1.076485613853039
$block2[16,5] = 1.069970665110734
$block2[17,0] = 1.055222873322578
$block2[17,1] = 1.068525756631092
$block2[17,2] = 1.068231784175447
$block2[17,3] = 0.993264023964098
$block2[17,4] = 1.076566746055707
$block2[17,5] = 1.070043185767125
$block2[18,0] = 1.055061535466545
$block2[18,1] = 1.068201760068174
$block2[18,2] = 1.067936985751164
$block2[18,3] = 0.9929783193494215
$block2[18,4] = 1.076203804445014
$block2[18,5] = 1.069718729091925
$block2[19,0] = 1.054534368732302
$block2[19,1] = 1.067146171499033
$block2[19,2] = 1.066976242662275
$block2[19,3] = 0.9920501090198102
$block2[19,4] = 1.07502208844988
$block2[19,5] = 1.068661641465938
$block2[20,0] = 1.054200883525624
$block2[20,1] = 1.066480743960525
$block2[20,2] = 1.066370388191296
$block2[20,3] = 0.9914670000341481
$block2[20,4] = 1.074277732149916
$block2[20,5] = 1.067995268944
$block2[21,0] = 1.054377863819817
$block2[21,1] = 1.066833665432551
$block2[21,2] = 1.066691734438255
$block2[21,3] = 0.9917760702887611
$block2[21,4] = 1.074672459969699
$block2[21,5] = 1.068348691604996
$block2[22,0] = 1.055070339135969
$block2[22,1] = 1.06821942794517
$block2[22,2] = 1.067953062489319
$block2[22,3] = 0.9929938892766441
$block2[22,4] = 1.076223593187877
$block2[22,5] = 1.069736422059333
$block2[23,0] = 1.055865121437749
$block2[23,1] = 1.069820072821605
$block2[23,2] = 1.069409035940159
$block2[23,3] = 0.9944092447426416
$block2[23,4] = 1.078017767004798
$block2[23,5] = 1.071339340035089

$ws.Range("I2:N25").Value = $block2

Write-Output "Applied updates to vm_pu sheet (rows 2-25)."